$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-string header labels in row 1 so that duplicated labels
# (e.g. "SE immediate", "LL95 immediate") become distinct per-section
# labels (CIN2+/CIN3+/CANCER prefixed), and several "Unweighted ..."
# columns get an explicit "Unweighted" prefix.
$ws.Range("V1").Value2 = "CIN2+ SE immediate"
$ws.Range("W1").Value2 = "CIN2+ LL95 immediate"
$ws.Range("X1").Value2 = "CIN2+ UL95 immediate"
$ws.Range("Z1").Value2 = "CIN2+ SE 1-year"
$ws.Range("AA1").Value2 = "CIN2+ LL95 1-year"
$ws.Range("AB1").Value2 = "CIN2+ UL95 1-year"
$ws.Range("AD1").Value2 = "CIN2+ SE 2-year"
$ws.Range("AE1").Value2 = "CIN2+ LL95 2-year"
$ws.Range("AF1").Value2 = "CIN2+ UL95 2-year"
$ws.Range("AH1").Value2 = "CIN2+ SE 3-year"
$ws.Range("AI1").Value2 = "CIN2+ LL95 3-year"
$ws.Range("AJ1").Value2 = "CIN2+ UL95 3-year"
$ws.Range("AL1").Value2 = "SCIN2+ E 4-year"
$ws.Range("AM1").Value2 = "CIN2+ LL95 4-year"
$ws.Range("AN1").Value2 = "CIN2+ UL95 4-year"
$ws.Range("AP1").Value2 = "CIN2+ SE 5-year"
$ws.Range("AQ1").Value2 = "CIN2+ LL95 5-year"
$ws.Range("AR1").Value2 = "CIN2+ UL95 5-year"
$ws.Range("AT1").Value2 = "CIN3+ SE immediate"
$ws.Range("AU1").Value2 = "CIN3+ LL95 immediate"
$ws.Range("AV1").Value2 = "CIN3+ UL95 immediate"
$ws.Range("AX1").Value2 = "CIN3+ SE 1-year"
$ws.Range("AY1").Value2 = "CIN3+ LL95 1-year"
$ws.Range("AZ1").Value2 = "CIN3+ UL95 1-year"
$ws.Range("BB1").Value2 = "CIN3+ SE 2-year"
$ws.Range("BC1").Value2 = "CIN3+ LL95 2-year"
$ws.Range("BD1").Value2 = "CIN3+ UL95 2-year"
$ws.Range("BF1").Value2 = "CIN3+ SE 3-year"
$ws.Range("BG1").Value2 = "CIN3+ L95 3-year"
$ws.Range("BH1").Value2 = "CIN3+ UL95 3-year"
$ws.Range("BJ1").Value2 = "CIN3+ SE 4-year"
$ws.Range("BK1").Value2 = "CIN3+ LL95 4-year"
$ws.Range("BL1").Value2 = "CIN3+ UL95 4-year"
$ws.Range("BN1").Value2 = "CIN3+ SE 5-year"
$ws.Range("BO1").Value2 = "CIN3+ LL95 5-year"
$ws.Range("BP1").Value2 = "CIN3+ UL95 5-year"
$ws.Range("BR1").Value2 = "CANCER SE immediate"
$ws.Range("BS1").Value2 = "CANCER LL95 immediate"
$ws.Range("BT1").Value2 = "CANCER UL95 immediate"
$ws.Range("BV1").Value2 = "CANCER SE 1-year"
$ws.Range("BW1").Value2 = "CANCER LL95 1-year"
$ws.Range("BX1").Value2 = "CANCER UL95 1-year"
$ws.Range("BZ1").Value2 = "CANCER SE 2-year"
$ws.Range("CA1").Value2 = "CANCER LL95 2-year"
$ws.Range("CB1").Value2 = "CANCER UL95 2-year"
$ws.Range("CD1").Value2 = "CANCER SE 3-year"
$ws.Range("CE1").Value2 = "CANCER LL95 3-year"
$ws.Range("CF1").Value2 = "CANCER UL95 3-year"
$ws.Range("CH1").Value2 = "CANCER SE 4-year"
$ws.Range("CI1").Value2 = "CANCER LL95 4-year"
$ws.Range("CJ1").Value2 = "CANCER UL95 4-year"
$ws.Range("CL1").Value2 = "CANCER SE 5-year"
$ws.Range("CM1").Value2 = "CANCER LL95 5-year"
$ws.Range("CN1").Value2 = "CANCER UL95 5-year"
$ws.Range("CS1").Value2 = "Unweighted %"
$ws.Range("CT1").Value2 = "Unweighted Informative N"
$ws.Range("CV1").Value2 = "Unweighted CIN2+ Prevalence Cases"
$ws.Range("CW1").Value2 = "Unweighted CIN2+ Incidence Cases"
$ws.Range("CX1").Value2 = "Unweighted CIN2+ Unknown Cases"
$ws.Range("CY1").Value2 = "Unweighted Number of CIN3+ Cases"
$ws.Range("CZ1").Value2 = "Unweighted CIN3+ Prevalence Cases"
$ws.Range("DA1").Value2 = "Unweighted CIN3+ Incidence Cases"
$ws.Range("DB1").Value2 = "Unweighted CIN3+ Unknown Cases"
$ws.Range("DD1").Value2 = "Unweighted Cancer Prevalence Cases"
$ws.Range("DE1").Value2 = "Unweighted Cancer Incidence Cases"
$ws.Range("DF1").Value2 = "Unweighted Cancer Unknown Cases"
# Update the frozen-pane selection / scroll position on the sheet view
# so the active selection is F1:DF1 (this also recomputes the pane
# topLeftCell to F2).
$ws.Range("F1:DF1").Select()
